$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Accuracy value (row 2)
$ws.Range("B2").Value = 0.9701149425287354

# Move current row 4 (Specificity) down to row 5, carrying label + value + format
$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A5").Value = "Specificity"
$ws.Range("B5").Value = 0.9840499920124092

# Move current row 3 (Sensitivity) down to row 4, carrying label + value + format
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A4").Value = "Sensitivity"
$ws.Range("B4").Value = 0.9451013922518159

# Row 3 now becomes the new "Error Rate" entry (reuse existing formatted cell)
$ws.Range("A3").Value = "Error Rate"
$ws.Range("B3").Value = 0.02988505747126437

# Add new row 6 "Geometric Mean"
$ws.Range("A2").Copy()
$ws.Range("A6").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A6").Value = "Geometric Mean"
$ws.Range("B6").Value = 0.9643111674053102

$excel.CutCopyMode = 0
